# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text before writing, so values like "12.30"
# or "2.80" keep their trailing zero instead of being auto-coerced to
# numbers by Excel's smart-entry parsing.
function Set-PriceText($cellRef, $text) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "69.551.22"
$ws.Range("E2").Value = "  +0.56%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.776.61"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "0.999"
$ws.Range("E4").Value = "  -0.36%  "

# Row 5 - BNB
Set-PriceText "D5" "664.65"
$ws.Range("E5").Value = "  +5.83%  "

# Row 6 - Solana
Set-PriceText "D6" "166.37"
$ws.Range("E6").Value = "  +1.58%  "

# Row 7 - LidoStakedEther
Set-PriceText "D7" "3.776.75"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - USDC (price unchanged)
$ws.Range("E8").Value = "  -0.12%  "

# Row 9 - XRP
Set-PriceText "D9" "0.527"
$ws.Range("E9").Value = "  +1.55%  "

# Row 10 - Dogecoin (price unchanged)
$ws.Range("E10").Value = "  -0.43%  "

# Row 11 - Cardano (price unchanged)
$ws.Range("E11").Value = "  +1.42%  "

# Row 12 - Toncoin
Set-PriceText "D12" "6.96"
$ws.Range("E12").Value = "  +5.16%  "

# Row 13 - ShibaInu (price unchanged)
$ws.Range("E13").Value = "  -2.31%  "

# Row 14 - Avalanche
Set-PriceText "D14" "35.16"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-PriceText "D15" "4.409.64"
$ws.Range("E15").Value = "  -0.40%  "

# Row 16 - WrappedEther
Set-PriceText "D16" "3.783.94"
$ws.Range("E16").Value = "  -2.74%  "

# Row 17 - WrappedBTC
Set-PriceText "D17" "69.416.81"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18 - Chainlink
Set-PriceText "D18" "17.72"
$ws.Range("E18").Value = "  -1.25%  "

# Row 19 - TRON (price unchanged)
$ws.Range("E19").Value = "  +0.59%  "

# Row 20 - Polkadot
Set-PriceText "D20" "7.11"
$ws.Range("E20").Value = "  +0.60%  "

# Row 21 - BitcoinCash
Set-PriceText "D21" "470.21"
$ws.Range("E21").Value = "  +0.55%  "

# Row 22 - Uniswap
Set-PriceText "D22" "9.62"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23 - Polygon
Set-PriceText "D23" "0.711"
$ws.Range("E23").Value = "  +1.28%  "

# Row 24 - PEPE (price unchanged)
$ws.Range("E24").Value = "  -2.81%  "

# Row 25 - Litecoin
Set-PriceText "D25" "82.26"
$ws.Range("E25").Value = "  -1.09%  "

# Row 26 - InternetComputer(DFINITY)
Set-PriceText "D26" "12.30"
$ws.Range("E26").Value = "  +2.33%  "

# Row 27 - RenderToken
Set-PriceText "D27" "10.26"
$ws.Range("E27").Value = "  +2.62%  "

# Row 28 - Fetch.AI (price unchanged)
$ws.Range("E28").Value = "  -1.24%  "

# Row 29 - Dai (price unchanged)
$ws.Range("E29").Value = "  +0.07%  "

# Row 30 - WrappedeETH
Set-PriceText "D30" "3.925.91"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31 - PancakeSwap
Set-PriceText "D31" "2.80"
$ws.Range("E31").Value = "  +4.95%  "

# Row 32 - ImmutableX
Set-PriceText "D32" "2.28"
$ws.Range("E32").Value = "  +3.38%  "

# Row 33 - NEARProtocol
Set-PriceText "D33" "7.24"
$ws.Range("E33").Value = "  -0.09%  "

# Row 34 - EthereumClassic
Set-PriceText "D34" "28.84"
$ws.Range("E34").Value = "  -0.11%  "

# Row 35 - Kaspa
Set-PriceText "D35" "0.177"
$ws.Range("E35").Value = "  +17.56%  "

# Row 36 - Binance-PegBSC-USD (price unchanged)
$ws.Range("E36").Value = "  +0.04%  "

# Row 37 - RenzoRestakedETH
Set-PriceText "D37" "3.732.63"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38 - Aptos
Set-PriceText "D38" "8.92"
$ws.Range("E38").Value = "  -0.48%  "

# Row 39 - Hedera (price unchanged)
$ws.Range("E39").Value = "  -1.44%  "

# Row 40 - dogwifhat (price unchanged)
$ws.Range("E40").Value = "  -0.24%  "

# Row 41 - Filecoin (price unchanged)
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 - now Mantle (was FirstDigitalUSD)
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceText "D42" "0.961"
$ws.Range("E42").Value = "  -0.70%  "

# Row 43 - now FirstDigitalUSD (was Mantle)
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-PriceText "D43" "0.999"
$ws.Range("E43").Value = "  -0.22%  "

# Row 44 - USDe (unchanged)

# Row 45 - Arweave
Set-PriceText "D45" "45.74"
$ws.Range("E45").Value = "  +7.70%  "

# Row 46 - Stacks (price unchanged)
$ws.Range("E46").Value = "  +4.89%  "

# Row 47 - Monero
Set-PriceText "D47" "157.18"
$ws.Range("E47").Value = "  +1.34%  "

# Row 48 - OKB
Set-PriceText "D48" "47.85"
$ws.Range("E48").Value = "  +2.13%  "

# Row 49 - TheGraph
Set-PriceText "D49" "0.298"
$ws.Range("E49").Value = "  +0.13%  "

# Row 50 - ONDO
Set-PriceText "D50" "1.41"
$ws.Range("E50").Value = "  +1.75%  "

# Row 51 - Cosmos
Set-PriceText "D51" "8.45"
$ws.Range("E51").Value = "  +0.62%  "
